$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = "S1234567A"
$ws.Cells.Item(4, 3).Value = 1
$ws.Cells.Item(4, 4).Value = "hi"
$ws.Cells.Item(4, 6).Value = 45767.27149535879

# Row 5
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = "S1234567A"
$ws.Cells.Item(5, 3).Value = 4
$ws.Cells.Item(5, 4).Value = "fe,feof"
$ws.Cells.Item(5, 6).Value = 45767.27272072917

# Row 6
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = "S1234567A"
$ws.Cells.Item(6, 3).Value = 2
$ws.Cells.Item(6, 4).Value = "?"
$ws.Cells.Item(6, 6).Value = 45767.28304898148
